$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

# The tile rows that used to describe the "Electron eat" animation (tileset
# "electr") are being repurposed for a new "Infotron eat" animation
# (tileset "infotr"). Only columns B (name) and C (tileset) actually change;
# every other column on these rows is a formula that recalculates
# automatically from B/C/D.
for ($i = 0; $i -lt 7; $i++) {
    $row = 69 + $i
    $n = $i + 1
    $ws.Cells.Item($row, 2).Value = "InfotronEat$n"
    $ws.Cells.Item($row, 3).Value = "infotr"
}

# Update the view state left behind by the edit: the selection moved to C67
# (with the window scrolled so row 31 is at the top, per the saved file).
$ws.Range("C67").Select()
